$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Employees")

$ws.Range("A2").Value = "Tom"
$ws.Range("A3").Value = "Madam"
